# edit.ps1 - apply the documented change to documentacao-site-individual.docx
#
# Summary of the change:
#  1. The first of the (previously empty) paragraphs that followed the
#     "Muitas franquias..." paragraph in the "Justificativa" section now
#     contains a new paragraph of text about why the author chose the
#     PlayStation 2 theme.
#  2. The next two (previously empty) paragraphs are removed entirely.
#  3. The "_GoBack" bookmark moves from the very end of the document to
#     right before the "Escopo" heading.

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Locate the anchor paragraph ("...agora os jogadores tem acesso a
#    estes titulos.") and the first empty paragraph right after it.
# -----------------------------------------------------------------
$anchorRng = $d.Content
$anchorRng.Find.Execute("estes títulos.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$anchorPara = $anchorRng.Paragraphs(1)
$targetPara = $anchorPara.Next()

$pRange = $targetPara.Range
$pRange.Collapse(1) | Out-Null
$startPos = $pRange.Start

# The three text segments that make up the new paragraph (split the same
# way the source document splits them around "PlayStation").
$seg1 = "Escolhi este tema porque o "
$seg2 = "PlayStation"
$seg3 = " 2 esteve presente em grande parte da minha vida. Praticamente todos os sábados de tarde, meu pai e eu costumávamos jogar os meus jogos favoritos para passarmos um tempo juntos. Era uma maneira que tínhamos de nos comunicar, já que não sabíamos direito o que falar um para o outro."
$fullText = $seg1 + $seg2 + $seg3

# Borrow the fully-specified run formatting (Arial 12pt, cs Arial 12pt)
# from the neighbouring text by copying its FormattedText into the
# target (empty) paragraph, then overwrite the copied text with the
# real content - this keeps the run properties (rFonts/sz/szCs) intact.
$pRange.FormattedText = $anchorRng.FormattedText
$fixRange = $d.Range($startPos, $startPos + "estes títulos.".Length)
$fixRange.Text = $fullText

# Re-apply FormattedText per-segment so the paragraph is split into three
# runs (mirroring how "PlayStation" is isolated from its neighbours in
# the rest of the document).
$seg1Start = $startPos
$seg1End = $seg1Start + $seg1.Length
$seg2Start = $seg1End
$seg2End = $seg2Start + $seg2.Length
$seg3Start = $seg2End
$seg3End = $seg3Start + $seg3.Length

$r1 = $d.Range($seg1Start, $seg1End)
$r2 = $d.Range($seg2Start, $seg2End)
$r3 = $d.Range($seg3Start, $seg3End)
$r1.FormattedText = $r1.FormattedText
$r2.FormattedText = $r2.FormattedText
$r3.FormattedText = $r3.FormattedText

# -----------------------------------------------------------------
# 2) Delete the two empty paragraphs that used to follow the one we
#    just filled in.
# -----------------------------------------------------------------
$p1 = $targetPara.Next()
$p2 = $p1.Next()
$p1.Range.Delete() | Out-Null
$p2.Range.Delete() | Out-Null

# -----------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document to just
#    before the "Escopo" heading.
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$escopoRng = $d.Content
$escopoRng.Find.Execute("Escopo", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$escopoStart = $escopoRng.Duplicate
$escopoStart.Collapse(1) | Out-Null
$d.Bookmarks.Add("_GoBack", $escopoStart) | Out-Null
